$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (A1:D1) to short machine-friendly names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case Spanish connector words (de, del, la, las, el, los, y) in
# state/municipality names throughout the data rows.
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B10").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "San José De Gracia"
$ws.Range("B16").Value = "Playas De Rosarito"
$ws.Range("B40").Value = "Amatenango De La Frontera"
$ws.Range("B41").Value = "Amatenango Del Valle"
$ws.Range("B44").Value = "Bejucal De Ocampo"
$ws.Range("B46").Value = "Benemérito De Las Américas"
$ws.Range("B56").Value = "Chiapa De Corzo"
$ws.Range("B63").Value = "Comitán De Domínguez"
$ws.Range("B92").Value = "Marqués De Comillas"
$ws.Range("B93").Value = "Mazapa De Madero"
$ws.Range("B96").Value = "Montecristo De Guerrero"
$ws.Range("B100").Value = "Ocozocoautla De Espinosa"
$ws.Range("B112").Value = "Salto De Agua"
$ws.Range("B113").Value = "San Cristóbal De Las Casas"
$ws.Range("B160").Value = "Coyame Del Sotol"
$ws.Range("B171").Value = "Guadalupe Y Calvo"
$ws.Range("B174").Value = "Hidalgo Del Parral"
$ws.Range("B199").Value = "San Francisco De Borja"
$ws.Range("B200").Value = "San Francisco De Conchos"
$ws.Range("B201").Value = "San Francisco Del Oro"
$ws.Range("B209").Value = "Valle De Zaragoza"
$ws.Range("B231").Value = "San Juan De Sabinas"
$ws.Range("B246").Value = "Villa De Álvarez"
$ws.Range("A248").Value = "Ciudad De México"
$ws.Range("B252").Value = "Cuajimalpa De Morelos"
$ws.Range("B267").Value = "Coneto De Comonfort"
$ws.Range("B281").Value = "Nombre De Dios"
$ws.Range("B285").Value = "Pánuco De Coronado"
$ws.Range("B292").Value = "San Juan De Guadalupe"
$ws.Range("B293").Value = "San Juan Del Río"
$ws.Range("B294").Value = "San Luis Del Cordero"
$ws.Range("B295").Value = "San Pedro Del Gallo"
$ws.Range("A305").Value = "Estado De México"
$ws.Range("B305").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B308").Value = "Almoloya De Alquisiras"
$ws.Range("B309").Value = "Almoloya De Juárez"
$ws.Range("B310").Value = "Almoloya Del Río"
$ws.Range("B317").Value = "Atizapán De Zaragoza"
$ws.Range("B325").Value = "Chapa De Mota"
$ws.Range("B331").Value = "Coacalco De Berriozábal"
$ws.Range("B338").Value = "Ecatepec De Morelos"
$ws.Range("B346").Value = "Ixtapan De La Sal"
$ws.Range("B347").Value = "Ixtapan Del Oro"
$ws.Range("B364").Value = "Naucalpan De Juárez"
$ws.Range("B378").Value = "San Antonio La Isla"
$ws.Range("B379").Value = "San Felipe Del Progreso"
$ws.Range("B380").Value = "San Martín De Las Pirámides"
$ws.Range("B382").Value = "San Simón De Guerrero"
$ws.Range("B384").Value = "Soyaniquilpan De Juárez"
$ws.Range("B393").Value = "Tenango Del Aire"
$ws.Range("B394").Value = "Tenango Del Valle"
$ws.Range("B408").Value = "Tlalnepantla De Baz"
$ws.Range("B414").Value = "Valle De Bravo"
$ws.Range("B415").Value = "Valle De Chalco Solidaridad"
$ws.Range("B416").Value = "Villa De Allende"
$ws.Range("B417").Value = "Villa Del Carbón"
$ws.Range("B430").Value = "San Miguel De Allende"
$ws.Range("B431").Value = "Apaseo El Alto"
$ws.Range("B432").Value = "Apaseo El Grande"
$ws.Range("B440").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B444").Value = "Jaral Del Progreso"
$ws.Range("B452").Value = "Purísima Del Rincón"
$ws.Range("B456").Value = "San Diego De La Unión"
$ws.Range("B458").Value = "San Francisco Del Rincón"
$ws.Range("B460").Value = "San Luis De La Paz"
$ws.Range("B462").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B464").Value = "Silao De La Victoria"
$ws.Range("B469").Value = "Valle De Santiago"
$ws.Range("B475").Value = "Acapulco De Juárez"
$ws.Range("B478").Value = "Ajuchitlán Del Progreso"
$ws.Range("B479").Value = "Alcozauca De Guerrero"
$ws.Range("B483").Value = "Atenango Del Río"
$ws.Range("B484").Value = "Atlamajalcingo Del Monte"
$ws.Range("B486").Value = "Atoyac De Álvarez"
$ws.Range("B487").Value = "Ayutla De Los Libres"
$ws.Range("B490").Value = "Buenavista De Cuéllar"
$ws.Range("B491").Value = "Chilapa De Álvarez"
$ws.Range("B492").Value = "Chilpancingo De Los Bravo"
$ws.Range("B493").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B498").Value = "Coyuca De Benítez"
$ws.Range("B499").Value = "Coyuca De Catalán"
$ws.Range("B503").Value = "Cuetzala Del Progreso"
$ws.Range("B504").Value = "Cutzamala De Pinzón"
$ws.Range("B510").Value = "Huitzuco De Los Figueroa"
$ws.Range("B511").Value = "Iguala De La Independencia"
$ws.Range("B513").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B514").Value = "Zihuatanejo De Azueta"
$ws.Range("B516").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B519").Value = "Mártir De Cuilapan"
$ws.Range("B532").Value = "Taxco De Alarcón"
$ws.Range("B534").Value = "Técpan De Galeana"
$ws.Range("B536").Value = "Tepecoacuilco De Trujano"
$ws.Range("B538").Value = "Tixtla De Guerrero"
$ws.Range("B542").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B543").Value = "Tlapa De Comonfort"
$ws.Range("B555").Value = "Agua Blanca De Iturbide"
$ws.Range("B562").Value = "Atotonilco De Tula"
$ws.Range("B563").Value = "Atotonilco El Grande"
$ws.Range("B569").Value = "Cuautepec De Hinojosa"
$ws.Range("B575").Value = "Huasca De Ocampo"
$ws.Range("B579").Value = "Huejutla De Reyes"
$ws.Range("B582").Value = "Jacala De Ledezma"
$ws.Range("B589").Value = "Mineral De La Reforma"
$ws.Range("B590").Value = "Mineral Del Chico"
$ws.Range("B591").Value = "Mineral Del Monte"
$ws.Range("B592").Value = "Mixquiahuala De Juárez"
$ws.Range("B593").Value = "Molango De Escamilla"
$ws.Range("B595").Value = "Nopala De Villagrán"
$ws.Range("B596").Value = "Omitlán De Juárez"
$ws.Range("B597").Value = "Pachuca De Soto"
$ws.Range("B600").Value = "Progreso De Obregón"
$ws.Range("B606").Value = "Santiago De Anaya"
$ws.Range("B607").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B611").Value = "Tenango De Doria"
$ws.Range("B613").Value = "Tepehuacán De Guerrero"
$ws.Range("B614").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B617").Value = "Tezontepec De Aldama"
$ws.Range("B626").Value = "Tula De Allende"
$ws.Range("B627").Value = "Tulancingo De Bravo"
$ws.Range("B628").Value = "Villa De Tezontepec"
$ws.Range("B632").Value = "Zacualtipán De Ángeles"
$ws.Range("B633").Value = "Zapotlán De Juárez"
$ws.Range("B638").Value = "Acatlán De Juárez"
$ws.Range("B639").Value = "Ahualulco De Mercado"
$ws.Range("B644").Value = "Atemajac De Brizuela"
$ws.Range("B647").Value = "Atotonilco El Alto"
$ws.Range("B649").Value = "Autlán De Navarro"
$ws.Range("B655").Value = "Cañadas De Obregón"
$ws.Range("B662").Value = "Concepción De Buenos Aires"
$ws.Range("B663").Value = "Cuautitlán De García Barragán"
$ws.Range("B672").Value = "Encarnación De Díaz"
$ws.Range("B679").Value = "Huejuquilla El Alto"
$ws.Range("B680").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B681").Value = "Ixtlahuacán Del Río"
$ws.Range("B685").Value = "Jilotlán De Los Dolores"
$ws.Range("B691").Value = "La Manzanilla De La Paz"
$ws.Range("B692").Value = "Lagos De Moreno"
$ws.Range("B700").Value = "Ojuelos De Jalisco"
$ws.Range("B705").Value = "San Cristóbal De La Barranca"
$ws.Range("B706").Value = "San Diego De Alejandría"
$ws.Range("B708").Value = "San Juan De Los Lagos"
$ws.Range("B709").Value = "San Juanito De Escobedo"
$ws.Range("B712").Value = "San Martín De Bolaños"
$ws.Range("B714").Value = "San Miguel El Alto"
$ws.Range("B715").Value = "San Sebastián Del Oeste"
$ws.Range("B716").Value = "Santa María De Los Ángeles"
$ws.Range("B717").Value = "Santa María Del Oro"
$ws.Range("B720").Value = "Talpa De Allende"
$ws.Range("B721").Value = "Tamazula De Gordiano"
$ws.Range("B724").Value = "Techaluta De Montenegro"
$ws.Range("B728").Value = "Teocuitatlán De Corona"
$ws.Range("B729").Value = "Tepatitlán De Morelos"
$ws.Range("B732").Value = "Tizapán El Alto"
$ws.Range("B733").Value = "Tlajomulco De Zúñiga"
$ws.Range("B745").Value = "Unión De San Antonio"
$ws.Range("B746").Value = "Unión De Tula"
$ws.Range("B747").Value = "Valle De Guadalupe"
$ws.Range("B748").Value = "Valle De Juárez"
$ws.Range("B753").Value = "Yahualica De González Gallo"
$ws.Range("B754").Value = "Zacoalco De Torres"
$ws.Range("B757").Value = "Zapotitlán De Vadillo"
$ws.Range("B758").Value = "Zapotlán Del Rey"
$ws.Range("B759").Value = "Zapotlán El Grande"
$ws.Range("B785").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B787").Value = "Cojumatlán De Régules"
$ws.Range("B854").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B880").Value = "Coatlán Del Río"
$ws.Range("B888").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B892").Value = "Puente De Ixtla"
$ws.Range("B898").Value = "Tetela Del Volcán"
$ws.Range("B900").Value = "Tlaltizapán De Zapata"
$ws.Range("B908").Value = "Zacualpan De Amilpas"
$ws.Range("B912").Value = "Amatlán De Cañas"
$ws.Range("B913").Value = "Bahía De Banderas"
$ws.Range("B917").Value = "Ixtlán Del Río"
$ws.Range("B924").Value = "Santa María Del Oro"
$ws.Range("B939").Value = "Ciénega De Flores"
$ws.Range("B952").Value = "Mier Y Noriega"
$ws.Range("B957").Value = "San Nicolás De Los Garza"
$ws.Range("B964").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B972").Value = "Ayoquezco De Aldama"
$ws.Range("B976").Value = "Capulálpam De Méndez"
$ws.Range("B978").Value = "Chalcatongo De Hidalgo"
$ws.Range("B979").Value = "Ciénega De Zimatlán"
$ws.Range("B982").Value = "Coicoyán De Las Flores"
$ws.Range("B985").Value = "Constancia Del Rosario"
$ws.Range("B988").Value = "Cuilápam De Guerrero"
$ws.Range("B989").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B990").Value = "El Barrio De La Soledad"
$ws.Range("B992").Value = "Eloxochitlán De Flores Magón"
$ws.Range("B993").Value = "Fresnillo De Trujano"
$ws.Range("B994").Value = "Guadalupe De Ramírez"
$ws.Range("B996").Value = "Guelatao De Juárez"
$ws.Range("B997").Value = "Guevea De Humboldt"
$ws.Range("B998").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B999").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B1000").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B1002").Value = "Huautla De Jiménez"
$ws.Range("B1004").Value = "Ixtlán De Juárez"
$ws.Range("B1005").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B1019").Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Range("B1021").Value = "Mariscala De Juárez"
$ws.Range("B1022").Value = "Mártires De Tacubaya"
$ws.Range("B1024").Value = "Mazatlán Villa De Flores"
$ws.Range("B1026").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B1027").Value = "Mixistlán De La Reforma"
$ws.Range("B1031").Value = "Nejapa De Madero"
$ws.Range("B1033").Value = "Oaxaca De Juárez"
$ws.Range("B1034").Value = "Ocotlán De Morelos"
$ws.Range("B1035").Value = "Pinotepa De Don Luis"
$ws.Range("B1037").Value = "Putla Villa De Guerrero"
$ws.Range("B1038").Value = "Reforma De Pineda"
$ws.Range("B1040").Value = "Rojas De Cuauhtémoc"
$ws.Range("B1045").Value = "San Agustín De Las Juntas"
$ws.Range("B1066").Value = "San Antonino El Alto"
$ws.Range("B1069").Value = "San Antonio De La Cal"
$ws.Range("B1075").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B1090").Value = "San Dionisio Del Mar"
$ws.Range("B1094").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B1101").Value = "San Francisco Del Mar"
$ws.Range("B1126").Value = "San José Del Peñasco"
$ws.Range("B1127").Value = "San José Del Progreso"
$ws.Range("B1138").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B1152").Value = "San Juan De Los Cués"
$ws.Range("B1153").Value = "San Juan Del Estado"
$ws.Range("B1154").Value = "San Juan Del Río"
$ws.Range("B1201").Value = "San Mateo Del Mar"
$ws.Range("B1219").Value = "San Miguel Del Puerto"
$ws.Range("B1220").Value = "San Miguel Del Río"
$ws.Range("B1222").Value = "San Miguel El Grande"
$ws.Range("B1248").Value = "San Pablo Villa De Mitla"
$ws.Range("B1256").Value = "San Pedro El Alto"
$ws.Range("B1283").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1284").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B1285").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B1303").Value = "Santa Ana Del Valle"
$ws.Range("B1323").Value = "Santa Cruz De Bravo"
$ws.Range("B1328").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B1334").Value = "Santa Inés De Zaragoza"
$ws.Range("B1335").Value = "Santa Inés Del Monte"
$ws.Range("B1337").Value = "Santa Lucía Del Camino"
$ws.Range("B1351").Value = "Santa María Del Rosario"
$ws.Range("B1352").Value = "Santa María Del Tule"
$ws.Range("B1360").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1362").Value = "Santa María La Asunción"
$ws.Range("B1402").Value = "Santiago Del Río"
$ws.Range("B1444").Value = "Santo Domingo De Morelos"
$ws.Range("B1469").Value = "Sitio De Xitlapehua"
$ws.Range("B1471").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1472").Value = "Tanetze De Zaragoza"
$ws.Range("B1474").Value = "Tataltepec De Valdés"
$ws.Range("B1475").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B1476").Value = "Teotitlán De Flores Magón"
$ws.Range("B1477").Value = "Teotitlán Del Valle"
$ws.Range("B1479").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B1480").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1481").Value = "Tlacolula De Matamoros"
$ws.Range("B1483").Value = "Tlalixtac De Cabrera"
$ws.Range("B1484").Value = "Totontepec Villa De Morelos"
$ws.Range("B1488").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1489").Value = "Villa De Etla"
$ws.Range("B1490").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1491").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1492").Value = "Villa De Zaachila"
$ws.Range("B1495").Value = "Villa Sola De Vega"
$ws.Range("B1496").Value = "Villa Talea De Castro"
$ws.Range("B1497").Value = "Villa Tejúpam De La Unión"
$ws.Range("B1500").Value = "Yutanduchi De Guerrero"
$ws.Range("B1501").Value = "Zapotitlán Del Río"
$ws.Range("B1504").Value = "Zimatlán De Álvarez"
$ws.Range("B1530").Value = "Ayotoxco De Guerrero"
$ws.Range("B1535").Value = "Chalchicomula De Sesma"
$ws.Range("B1545").Value = "Chila De La Sal"
$ws.Range("B1556").Value = "Cuapiaxtla De Madero"
$ws.Range("B1560").Value = "Cuayuca De Andrade"
$ws.Range("B1561").Value = "Cuetzalan Del Progreso"
$ws.Range("B1577").Value = "Huehuetlán El Chico"
$ws.Range("B1578").Value = "Huehuetlán El Grande"
$ws.Range("B1583").Value = "Huitzilan De Serdán"
$ws.Range("B1585").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1589").Value = "Izúcar De Matamoros"
$ws.Range("B1600").Value = "Los Reyes De Juárez"
$ws.Range("B1601").Value = "Mazapiltepec De Juárez"
$ws.Range("B1614").Value = "Palmar De Bravo"
$ws.Range("B1624").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1641").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1645").Value = "San Salvador El Seco"
$ws.Range("B1646").Value = "San Salvador El Verde"
$ws.Range("B1655").Value = "Tecali De Herrera"
$ws.Range("B1663").Value = "Tepanco De López"
$ws.Range("B1664").Value = "Tepango De Rodríguez"
$ws.Range("B1665").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1671").Value = "Tepexi De Rodríguez"
$ws.Range("B1673").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1674").Value = "Tetela De Ocampo"
$ws.Range("B1675").Value = "Teteles De Avila Castillo"
$ws.Range("B1680").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1692").Value = "Totoltepec De Guerrero"
$ws.Range("B1694").Value = "Tuzamapan De Galeana"
$ws.Range("B1698").Value = "Xayacatlán De Bravo"
$ws.Range("B1704").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1712").Value = "Zapotitlán De Méndez"
$ws.Range("B1720").Value = "Amealco De Bonfil"
$ws.Range("B1722").Value = "Cadereyta De Montes"
$ws.Range("B1728").Value = "Jalpan De Serra"
$ws.Range("B1729").Value = "Landa De Matamoros"
$ws.Range("B1732").Value = "Pinal De Amoles"
$ws.Range("B1735").Value = "San Juan Del Río"
$ws.Range("B1749").Value = "Armadillo De Los Infante"
$ws.Range("B1750").Value = "Axtla De Terrazas"
$ws.Range("B1756").Value = "Ciudad Del Maíz"
$ws.Range("B1766").Value = "Mexquitic De Carmona"
$ws.Range("B1772").Value = "San Ciro De Acosta"
$ws.Range("B1778").Value = "Santa María Del Río"
$ws.Range("B1780").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1790").Value = "Villa De Arista"
$ws.Range("B1791").Value = "Villa De Arriaga"
$ws.Range("B1792").Value = "Villa De Guadalupe"
$ws.Range("B1793").Value = "Villa De La Paz"
$ws.Range("B1794").Value = "Villa De Ramos"
$ws.Range("B1795").Value = "Villa De Reyes"
$ws.Range("B1858").Value = "Nacozari De García"
$ws.Range("B1868").Value = "San Miguel De Horcasitas"
$ws.Range("B1869").Value = "San Pedro De La Cueva"
$ws.Range("B1887").Value = "Jalpa De Méndez"
$ws.Range("B1925").Value = "Soto La Marina"
$ws.Range("B1933").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1935").Value = "Amaxac De Guerrero"
$ws.Range("B1936").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1942").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1950").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1954").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1955").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1958").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1963").Value = "San Pablo Del Monte"
$ws.Range("B1964").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B1972").Value = "Tepetitla De Lardizábal"
$ws.Range("B1975").Value = "Tetla De La Solidaridad"
$ws.Range("B1987").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1997").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B2001").Value = "Amatlán De Los Reyes"
$ws.Range("B2012").Value = "Boca Del Río"
$ws.Range("B2014").Value = "Camarón De Tejeda"
$ws.Range("B2018").Value = "Castillo De Teayo"
$ws.Range("B2020").Value = "Cazones De Herrera"
$ws.Range("B2027").Value = "Chinampa De Gorostiza"
$ws.Range("B2041").Value = "Cosamaloapan De Carpio"
$ws.Range("B2042").Value = "Cosautlán De Carvajal"
$ws.Range("B2059").Value = "Hueyapan De Ocampo"
$ws.Range("B2060").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B2061").Value = "Ignacio De La Llave"
$ws.Range("B2065").Value = "Ixhuacán De Los Reyes"
$ws.Range("B2066").Value = "Ixhuatlán De Madero"
$ws.Range("B2067").Value = "Ixhuatlán Del Café"
$ws.Range("B2068").Value = "Ixhuatlán Del Sureste"
$ws.Range("B2080").Value = "Juchique De Ferrer"
$ws.Range("B2083").Value = "Landero Y Coss"
$ws.Range("B2086").Value = "Las Vigas De Ramírez"
$ws.Range("B2087").Value = "Lerdo De Tejada"
$ws.Range("B2091").Value = "Martínez De La Torre"
$ws.Range("B2094").Value = "Medellín De Bravo"
$ws.Range("B2098").Value = "Mixtla De Altamirano"
$ws.Range("B2100").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B2111").Value = "Ozuluama De Mascareñas"
$ws.Range("B2115").Value = "Paso De Ovejas"
$ws.Range("B2116").Value = "Paso Del Macho"
$ws.Range("B2120").Value = "Poza Rica De Hidalgo"
$ws.Range("B2129").Value = "Sayula De Alemán"
$ws.Range("B2133").Value = "Soledad De Doblado"
$ws.Range("B2141").Value = "Tatahuicapan De Juárez"
$ws.Range("B2175").Value = "Vega De Alatorre"
$ws.Range("B2187").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B2188").Value = "Zozocolco De Hidalgo"
$ws.Range("B2210").Value = "Dzilam De Bravo"
$ws.Range("B2277").Value = "Cañitas De Felipe Pescador"
$ws.Range("B2279").Value = "Concepción Del Oro"
$ws.Range("B2281").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B2291").Value = "Jiménez Del Teul"
$ws.Range("B2297").Value = "Mezquital Del Oro"
$ws.Range("B2302").Value = "Moyahua De Estrada"
$ws.Range("B2303").Value = "Nochistlán De Mejía"
$ws.Range("B2304").Value = "Noria De Ángeles"
$ws.Range("B2315").Value = "Teúl De González Ortega"
$ws.Range("B2316").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B2318").Value = "Trinidad García De La Cadena"
$ws.Range("B2321").Value = "Villa De Cos"

# Remove the trailing metadata/footer rows (sample size, source, author,
# date) that followed the last data total row.
$ws.Range("A2329:A2334").EntireRow.Delete() | Out-Null
